# testdata excel data update
#
# The DATA worksheet's resume-file hyperlink text values (column H) are
# updated from local "D:\...\xxx" paths to relative "./resume/xxx" paths,
# and the active selection on that sheet moves from B12 to H12.
#
# In the original workbook the updated column-H cells used cell style
# index 1 (no quote prefix); in the edited workbook they use style index 2
# (same font, but with a leading quote prefix - the same style already
# used by the neighbouring column F cells). We replicate that by copying
# the number/format (quote-prefix) formatting from a cell that already
# carries that style after updating the cell's value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "DATA" sheet is already the active/selected sheet

# 1. Update the resume file path values in column H.
$ws.Range("H10").Value = "./resume/K_Thrinath.docx"
$ws.Range("H11").Value = "./resume/Sandeep_Resume.pdf"
$ws.Range("H12").Value = "./resume/Uday_Resume.docx"
$ws.Range("H16").Value = "./resume/Naresh_Resume.pdf"

# 2. Give those same four cells the "quote prefix" cell style (style index 2)
#    that other similar cells (e.g. F11) already use, by copying formats
#    from such a cell onto each of them (done after setting the value, since
#    assigning .Value resets direct formatting).
$fmtSource = $ws.Range("F11")

$fmtSource.Copy() | Out-Null
$ws.Range("H10").PasteSpecial(-4122) | Out-Null   # -4122 = xlPasteFormats

$fmtSource.Copy() | Out-Null
$ws.Range("H11").PasteSpecial(-4122) | Out-Null

$fmtSource.Copy() | Out-Null
$ws.Range("H12").PasteSpecial(-4122) | Out-Null

$fmtSource.Copy() | Out-Null
$ws.Range("H16").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# 3. Move the active selection on the DATA sheet from B12 to H12.
$ws.Range("H12").Select() | Out-Null
